$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the Date and Price columns as text so Excel does not
# auto-convert the incoming strings into date serials / numbers.
$ws.Range("C8:D9").NumberFormat = "@"

$ws.Range("A8").Value = "Other"
$ws.Range("B8").Value = "douchebag"
$ws.Range("C8").Value = "2023-03-22"
$ws.Range("D8").Value = "4000.0"
$ws.Range("E8").Value = "Card"

$ws.Range("A9").Value = "Clothing"
$ws.Range("B9").Value = "macbook"
$ws.Range("C9").Value = "2023-03-22"
$ws.Range("D9").Value = "7000.0"
$ws.Range("E9").Value = "Card"

# Remove the temporary text number-format so the new cells keep the
# same (default) styling as the rest of the sheet.
$ws.Range("C8:D9").ClearFormats()
